$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'34.373.72"
$ws.Range("E2").Value = "  +0.69%  "

# Row 3
$ws.Range("D3").Value = "'1.788.00"
$ws.Range("E3").Value = "  +0.28%  "

# Row 4
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").Value = "'225.91"
$ws.Range("E5").Value = "  -0.11%  "

# Row 6
$ws.Range("E6").Value = "  +1.01%  "

# Row 7
$ws.Range("E7").Value = "  -0.13%  "

# Row 8
$ws.Range("D8").Value = "'32.57"
$ws.Range("E8").Value = "  +1.48%  "

# Row 9
$ws.Range("E9").Value = "  +0.54%  "

# Row 10
$ws.Range("E10").Value = "  +0.22%  "

# Row 11
$ws.Range("D11").Value = "'0.0944"
$ws.Range("E11").Value = "  -0.61%  "

# Row 12
$ws.Range("D12").Value = "'2.043.99"
$ws.Range("E12").Value = "  +0.20%  "

# Row 13
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'11.02"
$ws.Range("E13").Value = "  +0.39%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "'1.770.30"
$ws.Range("E14").Value = "  -0.85%  "

# Row 15
$ws.Range("D15").Value = "'0.632"
$ws.Range("E15").Value = "  +1.56%  "

# Row 16
$ws.Range("D16").Value = "'34.349.15"
$ws.Range("E16").Value = "  +0.63%  "

# Row 17
$ws.Range("E17").Value = "  +2.19%  "

# Row 18
$ws.Range("D18").Value = "'68.14"
$ws.Range("E18").Value = "  +0.82%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0794"
$ws.Range("E19").Value = "  +0.79%  "

# Row 20
$ws.Range("D20").Value = "'244.00"
$ws.Range("E20").Value = "  -0.61%  "

# Row 21
$ws.Range("D21").Value = "'11.17"
$ws.Range("E21").Value = "  +2.33%  "

# Row 22
$ws.Range("E22").Value = "  -0.14%  "

# Row 23
$ws.Range("E23").Value = "  +0.64%  "

# Row 24
$ws.Range("E24").Value = "  +1.29%  "

# Row 25
$ws.Range("D25").Value = "'165.23"
$ws.Range("E25").Value = "  +1.92%  "

# Row 26
$ws.Range("E26").Value = "  +1.87%  "

# Row 27
$ws.Range("E27").Value = "  +1.01%  "

# Row 28
$ws.Range("E28").Value = "  +0.93%  "

# Row 30
$ws.Range("D30").Value = "'3.97"
$ws.Range("E30").Value = "  +6.35%  "

# Row 31
$ws.Range("D31").Value = "'0.0523"
$ws.Range("E31").Value = "  +0.85%  "

# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'3.79"
$ws.Range("E32").Value = "  +2.16%  "

# Row 33
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "'1.23"
$ws.Range("E33").Value = "  -0.02%  "

# Row 34
$ws.Range("D34").Value = "'1.81"
$ws.Range("E34").Value = "  +0.66%  "

# Row 35
$ws.Range("D35").Value = "'2.60"
$ws.Range("E35").Value = "  +3.17%  "

# Row 36
$ws.Range("D36").Value = "'1.401.36"
$ws.Range("E36").Value = "  -3.12%  "

# Row 37
$ws.Range("D37").Value = "'0.674"
$ws.Range("E37").Value = "  +3.29%  "

# Row 38
$ws.Range("D38").Value = "'1.06"
$ws.Range("E38").Value = "  +1.98%  "

# Row 39
$ws.Range("E39").Value = "  -0.56%  "

# Row 40
$ws.Range("D40").Value = "'84.30"
$ws.Range("E40").Value = "  +2.34%  "

# Row 41
$ws.Range("E41").Value = "  +0.71%  "

# Row 42
$ws.Range("D42").Value = "'2.78"
$ws.Range("E42").Value = "  +2.46%  "

# Row 43
$ws.Range("D43").Value = "'0.933"
$ws.Range("E43").Value = "  +2.11%  "

# Row 44
$ws.Range("D44").Value = "'13.73"
$ws.Range("E44").Value = "  +0.30%  "

# Row 45
$ws.Range("D45").Value = "'0.0525"
$ws.Range("E45").Value = "  +1.22%  "

# Row 46
$ws.Range("E46").Value = "  +3.45%  "

# Row 47
$ws.Range("D47").Value = "'5.96"
$ws.Range("E47").Value = "  -1.74%  "

# Row 48
$ws.Range("D48").Value = "'1.945.55"
$ws.Range("E48").Value = "  +0.30%  "

# Row 49
$ws.Range("D49").Value = "'104.69"
$ws.Range("E49").Value = "  -0.12%  "

# Row 51
$ws.Range("E51").Value = "  -3.05%  "
